# Update the "Förändrad" (Changed) date column (column C) for all data rows
# from the old date serial 45172 (2023-09-03) to the new date serial 45175
# (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 375
$oldValue = 45172
$newValue = 45175

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
